# Karma demo slides - "removing reference to commonjs" edit.
#
# This reproduces the text-level content changes from the target commit:
#   - Slide 8  ("Config Options - files"): merge the three runs that spell
#     out "Karma will include files in a <script> tag" into a single run,
#     and drop the " or CommonJS" tail so the paragraph reads
#     "...if you're using RequireJS".
#   - Slide 9  ("Config Options - files"): merge the three runs that spell
#     out "arma will serve the files with its webserver" into a single run.
#   - Slide 12 ("Config Options - coverage"): merge the five runs that
#     spell out ": { type: 'text-summary' //'html' }" into a single run.
#   - Slide 14 ("Use with Grunt / Gulp"): merge the three runs that spell
#     out "Use with Grunt / Gulp" into a single run.
#
# Each merge is done by replacing the exact character span (using
# TextRange.Characters(start,length)) with the identical text it already
# renders as; PowerPoint's text engine collapses that span back down to a
# single run carrying the formatting of the span's first original run,
# which is exactly what the target OOXML shows. The CommonJS removal is a
# genuine deletion of the trailing " or CommonJS" characters.

$p = $ppt.ActivePresentation

# --- Slide 8: "Config Options - files" (RequireJS/CommonJS slide) -------
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(2)
$tr8 = $sh8.TextFrame.TextRange

# Merge "Karma will include " + "files " + "in a <script> tag" (3 runs -> 1)
$tr8.Characters(59, 42).Text = "Karma will include files in a <script> tag"

# Remove " or CommonJS" after "RequireJS" (leaves "...using RequireJS")
$tr8.Characters(141, 12).Text = ""

# --- Slide 9: "Config Options - files" (webserver slide) ----------------
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(2)
$tr9 = $sh9.TextFrame.TextRange

# Merge "arma will serve the " + "files " + "with its webserver" (3 runs -> 1)
$tr9.Characters(60, 44).Text = "arma will serve the files with its webserver"

# --- Slide 12: "Config Options - coverage" -------------------------------
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(2)
$tr12 = $sh12.TextFrame.TextRange

# Merge ": { type: 'text-summary' " + "//'" + "html" + "' " + "}" (5 runs -> 1)
$tr12.Characters(91, 35).Text = [char]0x003A + " { type: " + [char]0x2018 + "text-summary" + [char]0x2019 + " //" + [char]0x2018 + "html" + [char]0x2019 + " }"

# --- Slide 14: "Use with Grunt / Gulp" title ------------------------------
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(1)
$tr14 = $sh14.TextFrame.TextRange

# Merge "Use with Grunt " + "/ " + "Gulp" (3 runs -> 1)
$tr14.Characters(1, 21).Text = "Use with Grunt / Gulp"

Write-Host "done"
